$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2154.55
$ws.Range("I40").Value = 2303.3333
$ws.Range("J40").Value = 2032.8182
$ws.Range("K40").Value = 2303.3333
$ws.Range("L40").Value = 2032.8182
$ws.Range("M40").Value = -2128.3333
$ws.Range("N40").Value = -2382.8182

$ws.Range("H64").Value = 4894.1816
$ws.Range("I64").Value = 3939.0908
$ws.Range("J64").Value = 5849.273
$ws.Range("K64").Value = 3939.0908
$ws.Range("L64").Value = 5849.273
$ws.Range("M64").Value = -3691.0908
$ws.Range("N64").Value = -6345.273

$ws.Range("H67").Value = 4894.1816
$ws.Range("I67").Value = 3939.0908
$ws.Range("J67").Value = 5849.273
$ws.Range("K67").Value = 3939.0908
$ws.Range("L67").Value = 5849.273
$ws.Range("M67").Value = -3081.0908
$ws.Range("N67").Value = -7565.273

$ws.Range("H76").Value = 3006075.8
$ws.Range("I76").Value = 3971125.2
$ws.Range("K76").Value = 3971125.2
$ws.Range("M76").Value = -3970810.2

$ws.Range("H79").Value = 3006075.8
$ws.Range("I79").Value = 3971125.2
$ws.Range("K79").Value = 3971125.2
$ws.Range("M79").Value = -3970033.2

$ws.Range("H135").Value = 1440.9
$ws.Range("I135").Value = 1610.1333
$ws.Range("K135").Value = 14491.1997
$ws.Range("M135").Value = -11956.1997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 37955.816
$ws.Range("I2").Value = 40944.28
$ws.Range("J2").Value = 600
$ws.Range("K2").Value = 40944.28
$ws.Range("L2").Value = 600
$ws.Range("M2").Value = -40831.28
$ws.Range("N2").Value = -826

$ws.Range("H9").Value = 27500

$ws.Range("H20").Value = 27500

$ws.Range("H63").Value = 5191.05
$ws.Range("I63").Value = 4224.6924
$ws.Range("J63").Value = 6985.7144
$ws.Range("K63").Value = 4224.6924
$ws.Range("L63").Value = 6985.7144
$ws.Range("M63").Value = -3538.6924
$ws.Range("N63").Value = -8357.714400000001

$ws.Range("H66").Value = 5191.05
$ws.Range("I66").Value = 4224.6924
$ws.Range("J66").Value = 6985.7144
$ws.Range("K66").Value = 21123.462
$ws.Range("L66").Value = 34928.572
$ws.Range("M66").Value = -17691.462
$ws.Range("N66").Value = -41792.572

$ws.Range("H92").Value = 32500
$ws.Range("J92").Value = 32500
$ws.Range("L92").Value = 32500
$ws.Range("N92").Value = -37492

$ws.Range("H116").Value = 37955.816
$ws.Range("I116").Value = 40944.28
$ws.Range("J116").Value = 600
$ws.Range("K116").Value = 40944.28
$ws.Range("L116").Value = 600
$ws.Range("M116").Value = -38650.28
$ws.Range("N116").Value = -5188

$ws.Range("H122").Value = 17287.428
$ws.Range("I122").Value = 23002.4
$ws.Range("K122").Value = 69007.20000000001
$ws.Range("M122").Value = -66557.20000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 37955.816
$ws.Range("I3").Value = 40944.28
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 40944.28
$ws.Range("L3").Value = 600
$ws.Range("M3").Value = -40830.28
$ws.Range("N3").Value = -828

$ws.Range("H92").Value = 27500
$ws.Range("J92").Value = 27500
$ws.Range("L92").Value = 27500
$ws.Range("N92").Value = -32492

$ws.Range("H105").Value = 246936.39
$ws.Range("I105").Value = 2795.3572
$ws.Range("J105").Value = 772778.6
$ws.Range("K105").Value = 2795.3572
$ws.Range("L105").Value = 772778.6
$ws.Range("M105").Value = -1048.3572
$ws.Range("N105").Value = -776272.6

$ws.Range("H109").Value = 34500
$ws.Range("J109").Value = 34500
$ws.Range("L109").Value = 34500
$ws.Range("N109").Value = -37274

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 17094.875
$ws.Range("I62").Value = 19501.46
$ws.Range("J62").Value = 6666.3335
$ws.Range("K62").Value = 19501.46
$ws.Range("L62").Value = 6666.3335
$ws.Range("M62").Value = -18877.46
$ws.Range("N62").Value = -7914.3335

$ws.Range("H65").Value = 17094.875
$ws.Range("I65").Value = 19501.46
$ws.Range("J65").Value = 6666.3335
$ws.Range("K65").Value = 97507.29999999999
$ws.Range("L65").Value = 33331.6675
$ws.Range("M65").Value = -94387.29999999999
$ws.Range("N65").Value = -39571.6675

$ws.Range("H132").Value = 1275.6305
$ws.Range("I132").Value = 825.3415
$ws.Range("K132").Value = 2476.0245
$ws.Range("M132").Value = 53.97550000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 902
$ws.Range("J19").Value = 902
$ws.Range("L19").Value = 2706
$ws.Range("N19").Value = -3054

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H70").Value = 5162.1333
$ws.Range("I70").Value = 5386.522
$ws.Range("J70").Value = 4424.857
$ws.Range("K70").Value = 5386.522
$ws.Range("L70").Value = 4424.857
$ws.Range("M70").Value = -5116.522
$ws.Range("N70").Value = -4964.857

$ws.Range("H73").Value = 5162.1333
$ws.Range("I73").Value = 5386.522
$ws.Range("J73").Value = 4424.857
$ws.Range("K73").Value = 5386.522
$ws.Range("L73").Value = 4424.857
$ws.Range("M73").Value = -4450.522
$ws.Range("N73").Value = -6296.857

$ws.Range("H80").Value = 3011.1765
$ws.Range("I80").Value = 2866
$ws.Range("K80").Value = 2866
$ws.Range("M80").Value = -1868

$ws.Range("H83").Value = 3011.1765
$ws.Range("I83").Value = 2866
$ws.Range("K83").Value = 14330
$ws.Range("M83").Value = -9338

$ws.Range("H122").Value = 1997
$ws.Range("I122").Value = 1997
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5991
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3541
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H21").Value = 4998.5
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 4998.5
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 4998.5
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -5346.5

$ws.Range("H96").Value = 32666.666
$ws.Range("J96").Value = 32666.666
$ws.Range("L96").Value = 32666.666
$ws.Range("N96").Value = -38158.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 7950
$ws.Range("J25").Value = 7950
$ws.Range("L25").Value = 7950
$ws.Range("N25").Value = -8536
